$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Model Name"
$ws.Range("C1").Value = "Response"

# Data row 2
$ws.Range("A2").Value = "What are the formats for loading a text file?"
$ws.Range("B2").Value = "deepseek1.5"
$ws.Range("C2").Value = "GEO can load several different types of ASCII files, such as CWLAS (LAS), tab-delimited, space-delimited and comma-delimited. The preferred data formats to request from your vendors are wholly structured LAS format or compatible-XML format."

# Data row 3
$ws.Range("A3").Value = "How many tracks can you define in one ODF (Output Database Type File)?"
$ws.Range("B3").Value = "llama3.2:latest"
$ws.Range("C3").Value = "According to the GEO application documentation, there is no specific limit mentioned on the number of tracks that can be defined in one ODF. However, it's recommended to keep the track count reasonable for efficient data management and analysis."

# Column widths (subtract Excel's automatic 5/6 padding so the stored
# OOXML <col width="..."> lands exactly on the target integer values)
$ws.Columns.Item(1).ColumnWidth = 72 - 5/6
$ws.Columns.Item(2).ColumnWidth = 17 - 5/6
$ws.Columns.Item(3).ColumnWidth = 247 - 5/6

# Header formatting: bold, borders, centered horizontally, top vertically
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Page margins (PageSetup properties are in points; 72pt = 1 inch)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

$ws.Range("A1").Select()
